$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 7.723999999999999
$ws.Range("B6").Value = 6.726000000000001
$ws.Range("B7").Value = 5.281
$ws.Range("C7").Value = -13.048
$ws.Range("B8").Value = 6.477000000000001
$ws.Range("C11").Value = -12.767
$ws.Range("C12").Value = -10.809
$ws.Range("E12").Value = 17.556
$ws.Range("E13").Value = 16.632
$ws.Range("E14").Value = 17.05
$ws.Range("C15").Value = -13.331
$ws.Range("B16").Value = 5.619
$ws.Range("E16").Value = 16.596
$ws.Range("E19").Value = 16.54
$ws.Range("B20").Value = 8.475
$ws.Range("C20").Value = -12.207
$ws.Range("E20").Value = 16.176
$ws.Range("B21").Value = 9.02
$ws.Range("C21").Value = -12.12
$ws.Range("C22").Value = -12.692
$ws.Range("E22").Value = 16.773
$ws.Range("C23").Value = -12.887
$ws.Range("B28").Value = 5.831
$ws.Range("B29").Value = 5.21
$ws.Range("C29").Value = -11.358
$ws.Range("B30").Value = 5.665
$ws.Range("B32").Value = 6.974000000000001
$ws.Range("C34").Value = -12.715
$ws.Range("E36").Value = 16.423
$ws.Range("B40").Value = 9.203999999999999
$ws.Range("C42").Value = -11.999
$ws.Range("C43").Value = -13.543
$ws.Range("E43").Value = 16.991
$ws.Range("C44").Value = -13.452
$ws.Range("C45").Value = -13.339
$ws.Range("B46").Value = 4.944000000000001
$ws.Range("C46").Value = -13.91
$ws.Range("E46").Value = 16.618
$ws.Range("C50").Value = -13.539
$ws.Range("E50").Value = 16.487
$ws.Range("B51").Value = 5.415
$ws.Range("C51").Value = -12.093
$ws.Range("B52").Value = 5.837000000000001
$ws.Range("B57").Value = 5.915
$ws.Range("C57").Value = -14.252
$ws.Range("B59").Value = 5.552999999999999
$ws.Range("B62").Value = 5.901999999999999
$ws.Range("C65").Value = -12.45
$ws.Range("B66").Value = 4.961
$ws.Range("C66").Value = -10.897
$ws.Range("C67").Value = -11.354
$ws.Range("B73").Value = 6.968999999999999
$ws.Range("B74").Value = 8.944000000000001
$ws.Range("E76").Value = 16.796
$ws.Range("B77").Value = 6.233
$ws.Range("C79").Value = -12.461
$ws.Range("C84").Value = -13.643
$ws.Range("C87").Value = -13.866
$ws.Range("B92").Value = 6.047
$ws.Range("C92").Value = -10.955
$ws.Range("E95").Value = 17.272
$ws.Range("C97").Value = -11.926
$ws.Range("E97").Value = 16.795
$ws.Range("E99").Value = 16.802
$ws.Range("B100").Value = 6.298
